$wb = $excel.ActiveWorkbook

# Rename the "Include ValueSets" sheets to "Include ValueSet #N"
$renames = @{
    "Include ValueSets"   = "Include ValueSet #0"
    "Include ValueSets 2" = "Include ValueSet #1"
    "Include ValueSets 3" = "Include ValueSet #2"
    "Include ValueSets 4" = "Include ValueSet #3"
    "Include ValueSets 5" = "Include ValueSet #4"
}

foreach ($ws in $wb.Worksheets) {
    $oldName = $ws.Name
    if ($renames.ContainsKey($oldName)) {
        $ws.Name = $renames[$oldName]
    }
}

# Update the Date value on the Metadata sheet (row with "Date" label in column A)
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-09-13T14:28:16+00:00"
